$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.61535666666667
$ws.Range("H2").Value = 31.84607
$ws.Range("I2").Value = 0.1437966543677284
$ws.Range("J2").Value = 0.1437966543677284
$ws.Range("M2").Value = 0.4652636666666667
$ws.Range("N2").Value = 1.395791
$ws.Range("O2").Value = 0.02604271297411062
$ws.Range("P2").Value = 0.02604271297411062
$ws.Range("Q2").Value = 4.938939765707778
$ws.Range("R2").Value = 44.45045789137
$ws.Range("S2").Value = 0.003744854996336142
$ws.Range("T2").Value = 0.003744854996336141

# Row 3
$ws.Range("G3").Value = 10.61535666666667
$ws.Range("H3").Value = 31.84607
$ws.Range("I3").Value = 0.1437966543677284
$ws.Range("J3").Value = 0.1437966543677284
$ws.Range("O3").Value = 0.09971126509087273
$ws.Range("P3").Value = 0.09971126509087272
$ws.Range("Q3").Value = 18.91000882803222
$ws.Range("R3").Value = 170.19007945229
$ws.Range("S3").Value = 0.01433814632284117
$ws.Range("T3").Value = 0.01433814632284117

# Row 4
$ws.Range("G4").Value = 10.61535666666667
$ws.Range("H4").Value = 31.84607
$ws.Range("I4").Value = 0.1437966543677284
$ws.Range("J4").Value = 0.1437966543677284
$ws.Range("M4").Value = 15.618761
$ws.Range("N4").Value = 46.856283
$ws.Range("O4").Value = 0.8742460219350168
$ws.Range("P4").Value = 0.8742460219350167
$ws.Range("Q4").Value = 165.7987187064233
$ws.Range("R4").Value = 1492.18846835781
$ws.Range("S4").Value = 0.1257136530485511
$ws.Range("T4").Value = 0.1257136530485511

# Row 5
$ws.Range("I5").Value = 0.6785840820776819
$ws.Range("J5").Value = 0.6785840820776819
$ws.Range("M5").Value = 0.4652636666666667
$ws.Range("N5").Value = 1.395791
$ws.Range("O5").Value = 0.02604271297411062
$ws.Range("P5").Value = 0.02604271297411062
$ws.Range("Q5").Value = 23.30712019751922
$ws.Range("R5").Value = 209.764081777673
$ws.Range("S5").Value = 0.0176721704783494
$ws.Range("T5").Value = 0.01767217047834939

# Row 6
$ws.Range("I6").Value = 0.6785840820776819
$ws.Range("J6").Value = 0.6785840820776819
$ws.Range("O6").Value = 0.09971126509087273
$ws.Range("P6").Value = 0.09971126509087272
$ws.Range("Q6").Value = 89.23734031972677
$ws.Range("R6").Value = 803.1360628775409
$ws.Range("S6").Value = 0.06766247729449429
$ws.Range("T6").Value = 0.06766247729449427

# Row 7
$ws.Range("I7").Value = 0.6785840820776819
$ws.Range("J7").Value = 0.6785840820776819
$ws.Range("M7").Value = 15.618761
$ws.Range("N7").Value = 46.856283
$ws.Range("O7").Value = 0.8742460219350168
$ws.Range("P7").Value = 0.8742460219350167
$ws.Range("Q7").Value = 782.4129972825276
$ws.Range("R7").Value = 7041.716975542749
$ws.Range("S7").Value = 0.5932494343048383
$ws.Range("T7").Value = 0.5932494343048382

# Row 8
$ws.Range("G8").Value = 13.11220933333333
$ws.Range("H8").Value = 39.336628
$ws.Range("I8").Value = 0.1776192635545896
$ws.Range("J8").Value = 0.1776192635545896
$ws.Range("M8").Value = 0.4652636666666667
$ws.Range("N8").Value = 1.395791
$ws.Range("O8").Value = 0.02604271297411062
$ws.Range("P8").Value = 0.02604271297411062
$ws.Range("Q8").Value = 6.100634592527556
$ws.Range("R8").Value = 54.90571133274801
$ws.Range("S8").Value = 0.004625687499425084
$ws.Range("T8").Value = 0.004625687499425084

# Row 9
$ws.Range("G9").Value = 13.11220933333333
$ws.Range("H9").Value = 39.336628
$ws.Range("I9").Value = 0.1776192635545896
$ws.Range("J9").Value = 0.1776192635545896
$ws.Range("O9").Value = 0.09971126509087273
$ws.Range("P9").Value = 0.09971126509087272
$ws.Range("Q9").Value = 23.35785805736844
$ws.Range("R9").Value = 210.220722516316
$ws.Range("S9").Value = 0.01771064147353727
$ws.Range("T9").Value = 0.01771064147353727

# Row 10
$ws.Range("G10").Value = 13.11220933333333
$ws.Range("H10").Value = 39.336628
$ws.Range("I10").Value = 0.1776192635545896
$ws.Range("J10").Value = 0.1776192635545896
$ws.Range("M10").Value = 15.618761
$ws.Range("N10").Value = 46.856283
$ws.Range("O10").Value = 0.8742460219350168
$ws.Range("P10").Value = 0.8742460219350167
$ws.Range("Q10").Value = 204.7964637593027
$ws.Range("R10").Value = 1843.168173833724
$ws.Range("S10").Value = 0.1552829345816273
$ws.Range("T10").Value = 0.1552829345816272
